$d = $word.ActiveDocument

# Change 1: "...Toolkit 2024!" -> "...Toolkit 2025!"
$range1 = $d.Content
$range1.Find.Execute("Toolkit 2024!", $true, $false, $false, $false, $false,
                      $true, 1, $false, "Toolkit 2025!", 2)

# Change 2: merge "b" + "oost" runs into a single "boost" run
# (text-wise, "to boost productivity" already reads correctly; this
# normalizes the underlying run split so it matches the target XML)
$range2 = $d.Content
$range2.Find.Execute("to boost productivity", $true, $false, $false, $false, $false,
                      $true, 1, $false, "to boost productivity", 2)
